# Weekly fruit/vegetable price update:
# a new "Uva" (grape) price record is inserted as row 142, pushing the
# existing rows 142-203 down to 143-204 (dimension grows from T203 to T204).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 142, shifting rows 142:203 down to 143:204.
$ws.Rows("142").Insert()

# Populate the newly inserted row with the new data entry.
$ws.Range("A142").Value = 11
$ws.Range("B142").Value = "Vega Monumental Concepción"
$ws.Range("C142").Value = "Bíobío"
$ws.Range("D142").Value = 45006
$ws.Range("E142").Value = 8
$ws.Range("F142").Value = "Fruta"
$ws.Range("G142").Value = 100109
$ws.Range("H142").Value = "Uva"
$ws.Range("I142").Value = 100109001
$ws.Range("J142").Value = "Uva"
$ws.Range("K142").Value = "Red Globe"
$ws.Range("L142").Value = "Primera"
$ws.Range("M142").Value = 180
$ws.Range("N142").Value = 9500
$ws.Range("O142").Value = 10000
$ws.Range("P142").Value = 9722
$ws.Range("Q142").Value = "$/bandeja 18 kilos"
$ws.Range("R142").Value = "Región de O'Higgins"
$ws.Range("S142").Value = 540
$ws.Range("T142").Value = 18
